$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2023-10-09 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-10 Tuesday", 2) | Out-Null

$tbl = $d.Tables.Item(1)

# Direct cell updates for values that are not unique in the document
# (global Find/Replace would be ambiguous for these)
$cell = $tbl.Cell(5, 5)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "64-24="

$cell = $tbl.Cell(17, 2)
$r = $cell.Range
$r.End = $r.End - 1
$r.Text = "27+52="

# Global Find/Replace for all remaining (unique) cell values
$d.Content.Find.Execute("84-59=", $true, $false, $false, $false, $false, $true, 1, $false, "72-26=", 2) | Out-Null
$d.Content.Find.Execute("61+26=", $true, $false, $false, $false, $false, $true, 1, $false, "18+23=", 2) | Out-Null
$d.Content.Find.Execute("5-0=", $true, $false, $false, $false, $false, $true, 1, $false, "50+16=", 2) | Out-Null
$d.Content.Find.Execute("11+71=", $true, $false, $false, $false, $false, $true, 1, $false, "90-72=", 2) | Out-Null
$d.Content.Find.Execute("29+25=", $true, $false, $false, $false, $false, $true, 1, $false, "85-15=", 2) | Out-Null
$d.Content.Find.Execute("27+32=", $true, $false, $false, $false, $false, $true, 1, $false, "52+21=", 2) | Out-Null
$d.Content.Find.Execute("46+41=", $true, $false, $false, $false, $false, $true, 1, $false, "4+90=", 2) | Out-Null
$d.Content.Find.Execute("14+80=", $true, $false, $false, $false, $false, $true, 1, $false, "82-2=", 2) | Out-Null
$d.Content.Find.Execute("1+28=", $true, $false, $false, $false, $false, $true, 1, $false, "98-66=", 2) | Out-Null
$d.Content.Find.Execute("61-36=", $true, $false, $false, $false, $false, $true, 1, $false, "51+16=", 2) | Out-Null
$d.Content.Find.Execute("36-9=", $true, $false, $false, $false, $false, $true, 1, $false, "51+24=", 2) | Out-Null
$d.Content.Find.Execute("16+8=", $true, $false, $false, $false, $false, $true, 1, $false, "86+13=", 2) | Out-Null
$d.Content.Find.Execute("87-73=", $true, $false, $false, $false, $false, $true, 1, $false, "81-32=", 2) | Out-Null
$d.Content.Find.Execute("21+8=", $true, $false, $false, $false, $false, $true, 1, $false, "7+23=", 2) | Out-Null
$d.Content.Find.Execute("69+9=", $true, $false, $false, $false, $false, $true, 1, $false, "30+8=", 2) | Out-Null
$d.Content.Find.Execute("34+1=", $true, $false, $false, $false, $false, $true, 1, $false, "98-70=", 2) | Out-Null
$d.Content.Find.Execute("75-4=", $true, $false, $false, $false, $false, $true, 1, $false, "91-70=", 2) | Out-Null
$d.Content.Find.Execute("86-16=", $true, $false, $false, $false, $false, $true, 1, $false, "30+68=", 2) | Out-Null
$d.Content.Find.Execute("74-35=", $true, $false, $false, $false, $false, $true, 1, $false, "7+16=", 2) | Out-Null
$d.Content.Find.Execute("44-7=", $true, $false, $false, $false, $false, $true, 1, $false, "63-22=", 2) | Out-Null
$d.Content.Find.Execute("16-15=", $true, $false, $false, $false, $false, $true, 1, $false, "97-84=", 2) | Out-Null
$d.Content.Find.Execute("53-36=", $true, $false, $false, $false, $false, $true, 1, $false, "64+30=", 2) | Out-Null
$d.Content.Find.Execute("44-9=", $true, $false, $false, $false, $false, $true, 1, $false, "39+46=", 2) | Out-Null
$d.Content.Find.Execute("9-1=", $true, $false, $false, $false, $false, $true, 1, $false, "84-46=", 2) | Out-Null
$d.Content.Find.Execute("1+9=", $true, $false, $false, $false, $false, $true, 1, $false, "64-57=", 2) | Out-Null
$d.Content.Find.Execute("81+12=", $true, $false, $false, $false, $false, $true, 1, $false, "20+51=", 2) | Out-Null
$d.Content.Find.Execute("46+39=", $true, $false, $false, $false, $false, $true, 1, $false, "92-84=", 2) | Out-Null
$d.Content.Find.Execute("59+16=", $true, $false, $false, $false, $false, $true, 1, $false, "3+90=", 2) | Out-Null
$d.Content.Find.Execute("24+62=", $true, $false, $false, $false, $false, $true, 1, $false, "5+86=", 2) | Out-Null
$d.Content.Find.Execute("27+70=", $true, $false, $false, $false, $false, $true, 1, $false, "76-15=", 2) | Out-Null
$d.Content.Find.Execute("76-45=", $true, $false, $false, $false, $false, $true, 1, $false, "15+40=", 2) | Out-Null
$d.Content.Find.Execute("9+12=", $true, $false, $false, $false, $false, $true, 1, $false, "66+1=", 2) | Out-Null
$d.Content.Find.Execute("55+14=", $true, $false, $false, $false, $false, $true, 1, $false, "1+16=", 2) | Out-Null
$d.Content.Find.Execute("26+20=", $true, $false, $false, $false, $false, $true, 1, $false, "66-49=", 2) | Out-Null
$d.Content.Find.Execute("62-29=", $true, $false, $false, $false, $false, $true, 1, $false, "2+86=", 2) | Out-Null
$d.Content.Find.Execute("42-38=", $true, $false, $false, $false, $false, $true, 1, $false, "30+63=", 2) | Out-Null
$d.Content.Find.Execute("44+40=", $true, $false, $false, $false, $false, $true, 1, $false, "58-26=", 2) | Out-Null
$d.Content.Find.Execute("32+59=", $true, $false, $false, $false, $false, $true, 1, $false, "4+72=", 2) | Out-Null
$d.Content.Find.Execute("58-39=", $true, $false, $false, $false, $false, $true, 1, $false, "57+20=", 2) | Out-Null
$d.Content.Find.Execute("82-19=", $true, $false, $false, $false, $false, $true, 1, $false, "2+31=", 2) | Out-Null
$d.Content.Find.Execute("1+50=", $true, $false, $false, $false, $false, $true, 1, $false, "75-28=", 2) | Out-Null
$d.Content.Find.Execute("40-25=", $true, $false, $false, $false, $false, $true, 1, $false, "49-8=", 2) | Out-Null
$d.Content.Find.Execute("17+13=", $true, $false, $false, $false, $false, $true, 1, $false, "74-34=", 2) | Out-Null
$d.Content.Find.Execute("39+38=", $true, $false, $false, $false, $false, $true, 1, $false, "9+50=", 2) | Out-Null
$d.Content.Find.Execute("4+86=", $true, $false, $false, $false, $false, $true, 1, $false, "6+58=", 2) | Out-Null
$d.Content.Find.Execute("94-15=", $true, $false, $false, $false, $false, $true, 1, $false, "65-14=", 2) | Out-Null
$d.Content.Find.Execute("88-63=", $true, $false, $false, $false, $false, $true, 1, $false, "50-15=", 2) | Out-Null
$d.Content.Find.Execute("24+21=", $true, $false, $false, $false, $false, $true, 1, $false, "21+5=", 2) | Out-Null
$d.Content.Find.Execute("84-7=", $true, $false, $false, $false, $false, $true, 1, $false, "27-18=", 2) | Out-Null
$d.Content.Find.Execute("92-44=", $true, $false, $false, $false, $false, $true, 1, $false, "27+17=", 2) | Out-Null
$d.Content.Find.Execute("71-48=", $true, $false, $false, $false, $false, $true, 1, $false, "46+14=", 2) | Out-Null
$d.Content.Find.Execute("66-52=", $true, $false, $false, $false, $false, $true, 1, $false, "97-37=", 2) | Out-Null
$d.Content.Find.Execute("57+39=", $true, $false, $false, $false, $false, $true, 1, $false, "93-13=", 2) | Out-Null
$d.Content.Find.Execute("74-0=", $true, $false, $false, $false, $false, $true, 1, $false, "15-13=", 2) | Out-Null
$d.Content.Find.Execute("68-54=", $true, $false, $false, $false, $false, $true, 1, $false, "39+33=", 2) | Out-Null
$d.Content.Find.Execute("25+63=", $true, $false, $false, $false, $false, $true, 1, $false, "11+53=", 2) | Out-Null
$d.Content.Find.Execute("90-31=", $true, $false, $false, $false, $false, $true, 1, $false, "18+44=", 2) | Out-Null
$d.Content.Find.Execute("90-54=", $true, $false, $false, $false, $false, $true, 1, $false, "1+18=", 2) | Out-Null
$d.Content.Find.Execute("68-65=", $true, $false, $false, $false, $false, $true, 1, $false, "34-3=", 2) | Out-Null
$d.Content.Find.Execute("10+4=", $true, $false, $false, $false, $false, $true, 1, $false, "60+8=", 2) | Out-Null
$d.Content.Find.Execute("6+83=", $true, $false, $false, $false, $false, $true, 1, $false, "21+75=", 2) | Out-Null
$d.Content.Find.Execute("4+32=", $true, $false, $false, $false, $false, $true, 1, $false, "29+67=", 2) | Out-Null
$d.Content.Find.Execute("39+39=", $true, $false, $false, $false, $false, $true, 1, $false, "81-38=", 2) | Out-Null
$d.Content.Find.Execute("33-7=", $true, $false, $false, $false, $false, $true, 1, $false, "58-56=", 2) | Out-Null
$d.Content.Find.Execute("60-9=", $true, $false, $false, $false, $false, $true, 1, $false, "66-13=", 2) | Out-Null
$d.Content.Find.Execute("22+1=", $true, $false, $false, $false, $false, $true, 1, $false, "14+27=", 2) | Out-Null
$d.Content.Find.Execute("33+43=", $true, $false, $false, $false, $false, $true, 1, $false, "14+39=", 2) | Out-Null
$d.Content.Find.Execute("36-34=", $true, $false, $false, $false, $false, $true, 1, $false, "52-40=", 2) | Out-Null
$d.Content.Find.Execute("53+8=", $true, $false, $false, $false, $false, $true, 1, $false, "45+40=", 2) | Out-Null
$d.Content.Find.Execute("93-93=", $true, $false, $false, $false, $false, $true, 1, $false, "19+43=", 2) | Out-Null
$d.Content.Find.Execute("65+26=", $true, $false, $false, $false, $false, $true, 1, $false, "31+64=", 2) | Out-Null
$d.Content.Find.Execute("9+31=", $true, $false, $false, $false, $false, $true, 1, $false, "57+19=", 2) | Out-Null
$d.Content.Find.Execute("68+8=", $true, $false, $false, $false, $false, $true, 1, $false, "48-27=", 2) | Out-Null
$d.Content.Find.Execute("32+34=", $true, $false, $false, $false, $false, $true, 1, $false, "74-72=", 2) | Out-Null
$d.Content.Find.Execute("98-18=", $true, $false, $false, $false, $false, $true, 1, $false, "68+29=", 2) | Out-Null
$d.Content.Find.Execute("8+54=", $true, $false, $false, $false, $false, $true, 1, $false, "7+0=", 2) | Out-Null
$d.Content.Find.Execute("70-11=", $true, $false, $false, $false, $false, $true, 1, $false, "0+17=", 2) | Out-Null
$d.Content.Find.Execute("51+36=", $true, $false, $false, $false, $false, $true, 1, $false, "59+23=", 2) | Out-Null
$d.Content.Find.Execute("28-14=", $true, $false, $false, $false, $false, $true, 1, $false, "46-36=", 2) | Out-Null
$d.Content.Find.Execute("28+16=", $true, $false, $false, $false, $false, $true, 1, $false, "29+57=", 2) | Out-Null
$d.Content.Find.Execute("31-7=", $true, $false, $false, $false, $false, $true, 1, $false, "67-38=", 2) | Out-Null
$d.Content.Find.Execute("75-61=", $true, $false, $false, $false, $false, $true, 1, $false, "22+40=", 2) | Out-Null
$d.Content.Find.Execute("69+17=", $true, $false, $false, $false, $false, $true, 1, $false, "42+11=", 2) | Out-Null
$d.Content.Find.Execute("5+81=", $true, $false, $false, $false, $false, $true, 1, $false, "55+11=", 2) | Out-Null
$d.Content.Find.Execute("12+53=", $true, $false, $false, $false, $false, $true, 1, $false, "51-0=", 2) | Out-Null
$d.Content.Find.Execute("4+4=", $true, $false, $false, $false, $false, $true, 1, $false, "50+38=", 2) | Out-Null
$d.Content.Find.Execute("33+33=", $true, $false, $false, $false, $false, $true, 1, $false, "97-14=", 2) | Out-Null
$d.Content.Find.Execute("63-13=", $true, $false, $false, $false, $false, $true, 1, $false, "81-22=", 2) | Out-Null
$d.Content.Find.Execute("74-59=", $true, $false, $false, $false, $false, $true, 1, $false, "75+17=", 2) | Out-Null
$d.Content.Find.Execute("89+8=", $true, $false, $false, $false, $false, $true, 1, $false, "34+3=", 2) | Out-Null
$d.Content.Find.Execute("67-27=", $true, $false, $false, $false, $false, $true, 1, $false, "47+14=", 2) | Out-Null
$d.Content.Find.Execute("96-58=", $true, $false, $false, $false, $false, $true, 1, $false, "92-50=", 2) | Out-Null
$d.Content.Find.Execute("99-5=", $true, $false, $false, $false, $false, $true, 1, $false, "42+21=", 2) | Out-Null
$d.Content.Find.Execute("51-31=", $true, $false, $false, $false, $false, $true, 1, $false, "66-45=", 2) | Out-Null
$d.Content.Find.Execute("49+36=", $true, $false, $false, $false, $false, $true, 1, $false, "34+3=", 2) | Out-Null
$d.Content.Find.Execute("6+68=", $true, $false, $false, $false, $false, $true, 1, $false, "81-16=", 2) | Out-Null
$d.Content.Find.Execute("72-6=", $true, $false, $false, $false, $false, $true, 1, $false, "11-10=", 2) | Out-Null
$d.Content.Find.Execute("96-34=", $true, $false, $false, $false, $false, $true, 1, $false, "90-4=", 2) | Out-Null
